$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text = "0M"
$t.Cell(2, 1).Range.Text = "0M"
$t.Cell(3, 1).Range.Text = "0M"
$t.Cell(4, 1).Range.Text = "3623"

$t.Cell(6, 1).Range.Text = "0.52663"
$t.Cell(7, 1).Range.Text = "0.07280"
$t.Cell(8, 1).Range.Text = "0.01158"
$t.Cell(9, 1).Range.Text = "0.42539"
$t.Cell(10, 1).Range.Text = "0.45507"
$t.Cell(11, 1).Range.Text = "0.47324"
$t.Cell(12, 1).Range.Text = "9.19949"

$t.Cell(44, 1).Range.Text = "99.26"
$t.Cell(45, 1).Range.Text = "9.2"
$t.Cell(46, 1).Range.Text = "1237"
